# Commit: "Added the ability to merge the rel xls files to one UW Relationship
# Master Report"
#
# The MODEL sheet in this report previously carried a bunch of #REF! formula
# cells (left behind from cells/ranges that were deleted elsewhere, e.g. in
# the source workbook(s) this report used to pull from before the relationship
# xls files were merged into one). Clear those stale formulas/cached errors
# out now, while keeping each cell's existing style/number format intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every address below held a formula such as "=#REF!" (or "=I6+I9" / "=K6+I9",
# which themselves resolved to #REF! via a broken downstream reference) plus
# a cached #REF! error value. ClearContents() drops the formula/value but
# keeps the cell's style (s="...") attribute, matching the diff which turns
# e.g. `<c r="B2" s="4" t="e"><f>#REF!</f><v>#REF!</v></c>` into
# `<c r="B2" s="4"/>`.
#
# A handful of these addresses are the top-left anchor of a merged range
# (B3:C3, B5:C5, B7:C7, D7:E7, B8:F8, E10:F10, B11:C11, E11:F14, B12:C12);
# clearing just the anchor cell leaves the merge's cached content alone, so
# for those we address the whole merged range instead.
$refCells = @(
    "B2", "E2", "G2",
    "B3:C3", "E3", "G3", "I3", "K3", "M3", "N3", "O3", "P3",
    "B5:C5", "G5",
    "B6", "G6", "I6", "J6", "K6", "M6", "N6", "O6", "P6",
    "B7:C7", "D7:E7", "F7",
    "B8:F8",
    "I9", "M9", "N9", "O9", "P9",
    "B10", "E10:F10",
    "B11:C11", "E11:F14",
    "B12:C12", "M12", "N12", "O12", "P12",
    "B13", "I13", "J13", "K13", "L13"
)

foreach ($addr in $refCells) {
    $ws.Range($addr).ClearContents()
}
